$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2-5) - analysis results with a new "Noun" method column
$rows = @(
    @{ A = 42606.574780092589; B = -52; C = 38; D = 60; E = 0;  F = 100; G = 3910; H = 1204; I = 215; J = 16; K = 25; L = 0; M = 4 },
    @{ A = 42606.575844907406; B = -26; C = 53; D = 46; E = 0;  F = 100; G = 3957; H = 1587; I = 308; J = 29; K = 25; L = 0; M = 4 },
    @{ A = 42606.580914351849; B = -26; C = 53; D = 46; E = 0;  F = 100; G = 4308; H = 1587; I = 308; J = 29; K = 25; L = 0; M = 4 },
    @{ A = 42606.581307870372; B = 2;   C = 51; D = 48; E = 50; F = 50;  G = 4848; H = 1905; I = 368; J = 35; K = 33; L = 1; M = 1 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = "Noun"
    $r++
}

# Column A now holds date/time values -> format + widen the column accordingly
$ws.Range("A1:A5").NumberFormat = "m/d/yy h:mm"
$ws.Columns.Item(1).ColumnWidth = 14
